$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 779848.3
$ws.Range("I132").Value = 2143.9333
$ws.Range("J132").Value = 2724109.2
$ws.Range("K132").Value = 6431.7999
$ws.Range("L132").Value = 8172327.600000001
$ws.Range("M132").Value = -3901.7999
$ws.Range("N132").Value = -8177387.600000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 26207.904
$ws.Range("I135").Value = 39028.37
$ws.Range("J135").Value = 3131.0667
$ws.Range("K135").Value = 351255.33
$ws.Range("L135").Value = 28179.6003
$ws.Range("M135").Value = -348720.33
$ws.Range("N135").Value = -33249.6003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2565848.8
$ws.Range("I137").Value = 3704925.8
$ws.Range("K137").Value = 11114777.4
$ws.Range("M137").Value = -11112227.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2194971
$ws.Range("I138").Value = 1551.878
$ws.Range("J138").Value = 4764405
$ws.Range("K138").Value = 4655.634
$ws.Range("L138").Value = 14293215
$ws.Range("M138").Value = 484.366
$ws.Range("N138").Value = -14303495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1522.54
$ws.Range("I32").Value = 1024.8193
$ws.Range("J32").Value = 3952.5881
$ws.Range("K32").Value = 1024.8193
$ws.Range("L32").Value = 3952.5881
$ws.Range("M32").Value = -737.8193000000001
$ws.Range("N32").Value = -4526.5881

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20041186
$ws.Range("I61").Value = 25026064
$ws.Range("J61").Value = 101672.8
$ws.Range("K61").Value = 25026064
$ws.Range("L61").Value = 101672.8
$ws.Range("M61").Value = -25025852
$ws.Range("N61").Value = -102096.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7799823
$ws.Range("I74").Value = 9288126
$ws.Range("K74").Value = 9288126
$ws.Range("M74").Value = -9287252

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7799823
$ws.Range("I77").Value = 9288126
$ws.Range("K77").Value = 46440630
$ws.Range("M77").Value = -46436262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2925909
$ws.Range("I122").Value = 1979.1666
$ws.Range("J122").Value = 13890646
$ws.Range("K122").Value = 5937.4998
$ws.Range("L122").Value = 41671938
$ws.Range("M122").Value = -3487.4998
$ws.Range("N122").Value = -41676838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 20041186
$ws.Range("I136").Value = 25026064
$ws.Range("J136").Value = 101672.8
$ws.Range("K136").Value = 75078192
$ws.Range("L136").Value = 305018.4
$ws.Range("M136").Value = -75075642
$ws.Range("N136").Value = -310118.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1738.449
$ws.Range("I134").Value = 1170.8529
$ws.Range("J134").Value = 3025
$ws.Range("K134").Value = 3512.5587
$ws.Range("L134").Value = 9075
$ws.Range("M134").Value = -977.5587000000005
$ws.Range("N134").Value = -14145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2566.75
$ws.Range("I31").Value = 1300.091
$ws.Range("K31").Value = 1300.091
$ws.Range("M31").Value = -1005.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2566.75
$ws.Range("I34").Value = 1300.091
$ws.Range("K34").Value = 1300.091
$ws.Range("M34").Value = -1098.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 6400
$ws.Range("J36").Value = 8000
$ws.Range("L36").Value = 8000
$ws.Range("N36").Value = -8776

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 6400
$ws.Range("J40").Value = 8000
$ws.Range("L40").Value = 8000
$ws.Range("N40").Value = -8320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 10000
$ws.Range("J42").Value = 10000
$ws.Range("L42").Value = 10000
$ws.Range("N42").Value = -11186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 9363.9
$ws.Range("J44").Value = 9363.9
$ws.Range("L44").Value = 9363.9
$ws.Range("N44").Value = -10247.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 27779992
$ws.Range("I58").Value = 50002572
$ws.Range("J58").Value = 1769.25
$ws.Range("K58").Value = 50002572
$ws.Range("L58").Value = 1769.25
$ws.Range("M58").Value = -50002369
$ws.Range("N58").Value = -2175.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3581.2
$ws.Range("I99").Value = 1768.6666
$ws.Range("J99").Value = 6300
$ws.Range("K99").Value = 1768.6666
$ws.Range("L99").Value = 6300
$ws.Range("M99").Value = -270.6666
$ws.Range("N99").Value = -9296

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3581.2
$ws.Range("I126").Value = 1768.6666
$ws.Range("J126").Value = 6300
$ws.Range("K126").Value = 5305.9998
$ws.Range("L126").Value = 18900
$ws.Range("M126").Value = -2835.9998
$ws.Range("N126").Value = -23840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 25230.326
$ws.Range("I132").Value = 1799.7142
$ws.Range("K132").Value = 5399.142599999999
$ws.Range("M132").Value = -2869.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 29238.436
$ws.Range("I134").Value = 1665.909
$ws.Range("J134").Value = 64920.53
$ws.Range("K134").Value = 4997.727000000001
$ws.Range("L134").Value = 194761.59
$ws.Range("M134").Value = -2462.727000000001
$ws.Range("N134").Value = -199831.59

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 27779992
$ws.Range("I136").Value = 50002572
$ws.Range("J136").Value = 1769.25
$ws.Range("K136").Value = 150007716
$ws.Range("L136").Value = 5307.75
$ws.Range("M136").Value = -150005166
$ws.Range("N136").Value = -10407.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1020.1429
$ws.Range("J131").Value = 1111.1476
$ws.Range("L131").Value = 3333.4428
$ws.Range("N131").Value = -13413.4428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 255.45454
$ws.Range("I107").Value = 125.55556
$ws.Range("K107").Value = 125.55556
$ws.Range("M107").Value = 1794.44444

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 70855.31
$ws.Range("I132").Value = 51589.6
$ws.Range("J132").Value = 113668
$ws.Range("K132").Value = 154768.8
$ws.Range("L132").Value = 341004
$ws.Range("M132").Value = -152238.8
$ws.Range("N132").Value = -346064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18796.396
$ws.Range("I132").Value = 1319.5435
$ws.Range("J132").Value = 85791
$ws.Range("K132").Value = 3958.6305
$ws.Range("L132").Value = 257373
$ws.Range("M132").Value = -1428.6305
$ws.Range("N132").Value = -262433

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1596.102
$ws.Range("I122").Value = 1076.5625
$ws.Range("J122").Value = 2574.0588
$ws.Range("K122").Value = 3229.6875
$ws.Range("L122").Value = 7722.176399999999
$ws.Range("M122").Value = -779.6875
$ws.Range("N122").Value = -12622.1764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 56839.316
$ws.Range("I132").Value = 56421.45
$ws.Range("J132").Value = 57303.61
$ws.Range("K132").Value = 169264.35
$ws.Range("L132").Value = 171910.83
$ws.Range("M132").Value = -166734.35
$ws.Range("N132").Value = -176970.83

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 30290.457
$ws.Range("I136").Value = 23621.523
$ws.Range("J136").Value = 41576.348
$ws.Range("K136").Value = 70864.569
$ws.Range("L136").Value = 124729.044
$ws.Range("M136").Value = -68314.569
$ws.Range("N136").Value = -129829.044
